$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: material + salary titles (Oxak, кг, 1, суммы, 1, 1, Shakhzod)
# Column A carries the same numeric/bordered style as the rows above it,
# so clone that formatting from A4 first, then set the number.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 3

$ws.Range("B5").Value = "Oxak"
$ws.Range("C5").Value = "кг"

# D5/F5/G5 hold the text "1" (not a number) — mirror an existing text "1"
# cell (F3) via a values-only paste so it round-trips as a string.
$ws.Range("F3").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("F5").PasteSpecial(-4163)
$ws.Range("G5").PasteSpecial(-4163)

$ws.Range("E5").Value = "суммы"
$ws.Range("H5").Value = "Shakhzod"
